$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 139
$ws.Range("I2").Value = 139
$ws.Range("K2").Value = 139
$ws.Range("M2").Value = -26
$ws.Range("H17").Value = 2089.3635
$ws.Range("J17").Value = 2089.3635
$ws.Range("L17").Value = 6268.0905
$ws.Range("N17").Value = -6604.0905
$ws.Range("H53").Value = 23811050
$ws.Range("J53").Value = 1782.5555
$ws.Range("L53").Value = 1782.5555
$ws.Range("N53").Value = -3056.5555
$ws.Range("H92").Value = 169.58824
$ws.Range("I92").Value = 125.6
$ws.Range("K92").Value = 125.6
$ws.Range("M92").Value = 1122.4
$ws.Range("H96").Value = 798.1
$ws.Range("I96").Value = 810.5
$ws.Range("J96").Value = 748.5
$ws.Range("K96").Value = 2431.5
$ws.Range("L96").Value = 2245.5
$ws.Range("M96").Value = -1058.5
$ws.Range("N96").Value = -4991.5
$ws.Range("H107").Value = 112310.22
$ws.Range("I107").Value = 126299
$ws.Range("K107").Value = 126299
$ws.Range("M107").Value = -124379
$ws.Range("H132").Value = 1571.1455
$ws.Range("I132").Value = 1571.1455
$ws.Range("K132").Value = 4713.4365
$ws.Range("M132").Value = -2183.4365
$ws.Range("H133").Value = 49999.918
$ws.Range("J133").Value = 49999.918
$ws.Range("L133").Value = 49999.918
$ws.Range("N133").Value = -60119.918
$ws.Range("H138").Value = 5595.9785
$ws.Range("I138").Value = 2869.6
$ws.Range("J138").Value = 6894.254
$ws.Range("K138").Value = 8608.799999999999
$ws.Range("L138").Value = 20682.762
$ws.Range("M138").Value = -3468.799999999999
$ws.Range("N138").Value = -30962.762

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2325.8572
$ws.Range("I45").Value = 1927.6875
$ws.Range("K45").Value = 1927.6875
$ws.Range("M45").Value = -1550.6875
$ws.Range("H132").Value = 5960.7144
$ws.Range("I132").Value = 5065.3335
$ws.Range("K132").Value = 15196.0005
$ws.Range("M132").Value = -12666.0005

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 29764.053
$ws.Range("I134").Value = 3400.9714
$ws.Range("J134").Value = 337333.34
$ws.Range("K134").Value = 10202.9142
$ws.Range("L134").Value = 1012000.02
$ws.Range("M134").Value = -7667.914199999999
$ws.Range("N134").Value = -1017070.02

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 40012.555
$ws.Range("I59").Value = 30000
$ws.Range("J59").Value = 41264.125
$ws.Range("K59").Value = 30000
$ws.Range("L59").Value = 41264.125
$ws.Range("M59").Value = -28855
$ws.Range("N59").Value = -43554.125
$ws.Range("H141").Value = 400579.9
$ws.Range("J141").Value = 419447.22
$ws.Range("L141").Value = 419447.22
$ws.Range("N141").Value = -429807.22

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H42").Value = 5133.3335
$ws.Range("J42").Value = 7500
$ws.Range("L42").Value = 22500
$ws.Range("N42").Value = -23568
$ws.Range("H74").Value = 6015
$ws.Range("J74").Value = 6015
$ws.Range("L74").Value = 18045
$ws.Range("N74").Value = -20167
$ws.Range("H77").Value = 6015
$ws.Range("J77").Value = 6015
$ws.Range("L77").Value = 54135
$ws.Range("N77").Value = -64743
$ws.Range("H86").Value = 1251.5
$ws.Range("J86").Value = 1251.5
$ws.Range("L86").Value = 3754.5
$ws.Range("N86").Value = -6126.5
$ws.Range("H89").Value = 1251.5
$ws.Range("J89").Value = 1251.5
$ws.Range("L89").Value = 11263.5
$ws.Range("N89").Value = -23119.5
$ws.Range("H93").Value = 1826.5
$ws.Range("I93").Value = 1826.5
$ws.Range("K93").Value = 5479.5
$ws.Range("M93").Value = -3607.5
$ws.Range("H96").Value = 291324.16
$ws.Range("J96").Value = 337916.34
$ws.Range("L96").Value = 1013749.02
$ws.Range("N96").Value = -1017867.02
$ws.Range("H100").Value = 3764
$ws.Range("J100").Value = 3764
$ws.Range("L100").Value = 11292
$ws.Range("N100").Value = -12914
$ws.Range("H101").Value = 5980.6665
$ws.Range("I101").Value = 4413
$ws.Range("K101").Value = 13239
$ws.Range("M101").Value = -10805
$ws.Range("H106").Value = 32163.334
$ws.Range("J106").Value = 34245
$ws.Range("L106").Value = 102735
$ws.Range("N106").Value = -104627
$ws.Range("H110").Value = 2528
$ws.Range("I110").Value = 2528
$ws.Range("K110").Value = 7584
$ws.Range("M110").Value = -3494
$ws.Range("H119").Value = 19676.334
$ws.Range("I119").Value = 19676.334
$ws.Range("K119").Value = 59029.00199999999
$ws.Range("M119").Value = -54191.00199999999
$ws.Range("H129").Value = 19667992
$ws.Range("J129").Value = 145193
$ws.Range("L129").Value = 435579
$ws.Range("N129").Value = -445579
$ws.Range("H138").Value = 2422.7778
$ws.Range("I138").Value = 2036
$ws.Range("K138").Value = 6108
$ws.Range("M138").Value = -968
$ws.Range("H139").Value = 6073.057
$ws.Range("I139").Value = 2516.3125
$ws.Range("K139").Value = 7548.9375
$ws.Range("M139").Value = -2408.9375

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 142876240
$ws.Range("I70").Value = 4374.5
$ws.Range("K70").Value = 4374.5
$ws.Range("M70").Value = -4104.5
$ws.Range("H73").Value = 142876240
$ws.Range("I73").Value = 4374.5
$ws.Range("K73").Value = 4374.5
$ws.Range("M73").Value = -3438.5
$ws.Range("H126").Value = 4494.6665
$ws.Range("I126").Value = 5000
$ws.Range("J126").Value = 4431.5
$ws.Range("K126").Value = 15000
$ws.Range("L126").Value = 13294.5
$ws.Range("M126").Value = -12530
$ws.Range("N126").Value = -18234.5
$ws.Range("H132").Value = 61659.277
$ws.Range("I132").Value = 6462.7646
$ws.Range("J132").Value = 1000000
$ws.Range("K132").Value = 19388.2938
$ws.Range("L132").Value = 3000000
$ws.Range("M132").Value = -16858.2938
$ws.Range("N132").Value = -3005060
$ws.Range("H134").Value = 65000
$ws.Range("J134").Value = 65000
$ws.Range("L134").Value = 195000
$ws.Range("N134").Value = -200070

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 485272.72
$ws.Range("I20").Value = 18250
$ws.Range("K20").Value = 18250
$ws.Range("M20").Value = -18024
$ws.Range("H46").Value = 1703.8235
$ws.Range("I46").Value = 1729.9
$ws.Range("K46").Value = 1729.9
$ws.Range("M46").Value = -1541.9
$ws.Range("H93").Value = 3164.2896
$ws.Range("I93").Value = 2869
$ws.Range("J93").Value = 3889.0908
$ws.Range("K93").Value = 2869
$ws.Range("L93").Value = 3889.0908
$ws.Range("M93").Value = -1621
$ws.Range("N93").Value = -6385.0908
$ws.Range("H100").Value = 4465.8667
$ws.Range("I100").Value = 2797.6
$ws.Range("K100").Value = 2797.6
$ws.Range("M100").Value = -2256.6

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 17500
$ws.Range("I45").Value = 17000
$ws.Range("J45").Value = 18000
$ws.Range("K45").Value = 17000
$ws.Range("L45").Value = 18000
$ws.Range("M45").Value = -16509
$ws.Range("N45").Value = -18982
$ws.Range("H100").Value = 691.3333
$ws.Range("I100").Value = 706.5
$ws.Range("J100").Value = 661
$ws.Range("K100").Value = 1413
$ws.Range("L100").Value = 1322
$ws.Range("M100").Value = -872
$ws.Range("N100").Value = -2404
$ws.Range("H132").Value = 26441.707
$ws.Range("I132").Value = 1338.7241
$ws.Range("K132").Value = 4016.1723
$ws.Range("M132").Value = -1486.1723
$ws.Range("H137").Value = 61249.75
$ws.Range("J137").Value = 61249.75
$ws.Range("L137").Value = 61249.75
$ws.Range("N137").Value = -71449.75
